$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume table refresh (scheduled data pull).
# Column D holds prices formatted as plain text (so values like "230.71"
# or "2.30" keep their trailing digits instead of being read as numbers),
# so force text format before writing any value that would otherwise be
# auto-detected as a number by Excel.

$ws.Range("D2").Value = '43.820.58'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '2.260.34'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.71'
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("E6").Value = '  +2.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.47'
$ws.Range("E7").Value = '  +5.35%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.451'
$ws.Range("E9").Value = '  +7.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0981'
$ws.Range("E10").Value = '  +5.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.03'
$ws.Range("E11").Value = '  -1.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.07'
$ws.Range("E12").Value = '  +15.21%  '

$ws.Range("E13").Value = '  +1.38%  '

$ws.Range("D14").Value = '2.597.30'
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.65'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.09'
$ws.Range("E16").Value = '  +5.29%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +3.60%  '

$ws.Range("D18").Value = '2.264.71'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").Value = '43.871.54'
$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").Value = '0.0₃0992'
$ws.Range("E20").Value = '  +6.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.41'
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("E22").Value = '  -2.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.30'
$ws.Range("E23").Value = '  -1.30%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("E25").Value = '  -3.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  +0.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  +2.56%  '

$ws.Range("E28").Value = '  +22.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.13'

$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.95'
$ws.Range("E31").Value = '  +2.33%  '

$ws.Range("E32").Value = '  -2.92%  '

$ws.Range("E33").Value = '  +3.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0706'
$ws.Range("E34").Value = '  +7.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.79'
$ws.Range("E35").Value = '  +0.79%  '

$ws.Range("E36").Value = '  -3.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.75'
$ws.Range("E37").Value = '  +4.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.49'
$ws.Range("E38").Value = '  +0.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.30'
$ws.Range("E39").Value = '  -3.63%  '

$ws.Range("E40").Value = '  +4.54%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000224'
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.37'
$ws.Range("E43").Value = '  +4.57%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0971'
$ws.Range("E44").Value = '  -1.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.22'
$ws.Range("E45").Value = '  -5.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '97.94'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.41'
$ws.Range("E48").Value = '  -2.43%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.38'
$ws.Range("E49").Value = '  +6.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.25'
$ws.Range("E50").Value = '  +7.20%  '

$ws.Range("D51").Value = '1.439.33'
$ws.Range("E51").Value = '  -2.26%  '
